$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update ingredient / nutrition text cells (column B = Ingredients, column F = LeaveEmpty) ---
# Row 7: Cheese Burger - strip bracketed sub-ingredient detail from Beef Patty description
$ws.Cells.Item(7, 2).Value = "Sesame Burger Bun / Beef Patty / Pickles / Cheddar Cheese / Red Onion"

# Row 8: Veggie Burger - strip bracketed sub-ingredient detail from Black Bean Patty description
$ws.Cells.Item(8, 2).Value = "Sesame Burger Bun / Spicy Black Bean Patty / Pickles / Cheddar Cheese / Red Onion"

# Row 9: Fries - simplified ingredient description
$ws.Cells.Item(9, 2).Value = "McCain Gold Crisp Fries: Potato / Oil / Seasoned Coating"

# Row 11: Spring Roll - strip bracketed sub-ingredient detail
$ws.Cells.Item(11, 2).Value = "Spring Roll / Sweet and Spicy Thai Sauce."

# Row 12: Tater Tot - reworded / recapitalized ingredient description
$ws.Cells.Item(12, 2).Value = "Potatoes / Vegetable Oil  / Salt / Dehydrated Onion / Corn Starch / Sodium Phosphate / Dextrose."

# Row 16: 5 Cheese Veggie Calzone - fix typo "Spinave" -> "Spinach"
$ws.Cells.Item(16, 2).Value = "Spinach / Tomato /Onions / Peppers / Tomato Sauce / Mozzarella / Feta / Parmesan"

# Row 15: Pepperoni & Sausage Calzone - replace placeholder LeaveEmpty value
$ws.Cells.Item(15, 6).Value = "Sausage_and_Pepperoni_Calzone"

# Row 16: 5 Cheese Veggie Calzone - replace placeholder LeaveEmpty value
$ws.Cells.Item(16, 6).Value = "5_Cheese_Veggie_Calzone"

# --- Update the active selection / view (was F17, now F16) ---
$ws.Range("F16").Select()
